# Email 2-1 [TEMPLATE] Partner email - reminder to RSVP.docx
# Traditional Chinese -> Simplified Chinese re-translation.
#
# Strategy: walk a single persistent Range ($r = $d.Content) forward through
# the document, issuing Find.Execute(..., Replace:=1 [wdReplaceOne]) calls in
# strict left-to-right document order. Because the range cursor only ever
# advances, duplicate substrings earlier in the doc are already consumed by
# the time a later call looks for the "same" text further on, so each call
# unambiguously hits the one run it is meant for.

$d = $word.ActiveDocument
$r = $d.Content

function Replace1($old, $new) {
    $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

Replace1 "英文" "英语"
Replace1 " / 葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語" " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语"
Replace1 "英文" "英语"
Replace1 "簡介" "简要"
Replace1 "向目標國的合作夥伴發送的電子郵件，以提醒他們提交RSVP。 將通過 customer.io 發送" "给未回复RSVP的目标国家合作伙伴的邮件，提醒他们发送RSVP。 将通过 customer.io 发送"
Replace1 "目標受眾" "目标受众"
Replace1 "尚未回覆的受邀合作夥伴" "尚未确认RSVP的受邀合作伙伴"
Replace1 "主題行" "主题行"
Replace1 ": 提醒：請回覆 " ": 提醒：请回复RSVP以参加"
Replace1 "[活動名稱]" "[EVENT NAME]"
Replace1 "不要拖延！ 立即預訂！" "不要拖延！ 立即预订！"
Replace1 "[合作夥伴姓名]" "[合作伙伴姓名]"
Replace1 "希望您和我們一樣為將於 " "希望您和我们一样为将于"
Replace1 "[活動名稱]" "[活动名称]"
Replace1 " 舉行的 " "举行的活动感到兴奋，活动定于"
Replace1 " 感到興奮！" "！"
Replace1 "希望您和我們一樣為將於 " "希望您和我们一样为将于"
Replace1 "[活動名稱]" "[活动名称]"
Replace1 " 從 " "举行的活动感到兴奋，活动定于"
Replace1 " 到 " "到"
Replace1 " 舉行的活動感到興奮！" "！"
Replace1 "請於 [" "请于["
Replace1 "DD Mmm YYYY" "日月年"
Replace1 "] 前確認您會參加這個備受期待的活動，因為名額有限，先到先得。" "]前确认您会参加这个备受期待的活动，因为名额有限，先到先得。"
Replace1 "立即回覆" "立即确认"
Replace1 "如有任何疑問，請與我們聯繫：" "如有任何疑问，请通过"
Replace1 "即時聊天" "实时聊天"
Replace1 " 或 " "或"
Replace1 "。" "联系我们。 "
Replace1 "如有任何疑問，請聯繫您的區域經理 " "如有任何疑问，请联系您的区域经理"
Replace1 "，電子郵件 " "，通过"
Replace1 "[電子郵件地址]" "[电子邮件地址]"
Replace1 " 或 " "或"
Replace1 "[WHATSAPP 號碼]" "[WHATSAPP 号码]"
Replace1 "期待在 " "期待在"
Replace1 "[活動名稱]" "[EVENT NAME]"
Replace1 " 與您見面！" "见到您！ "
Replace1 "期待在 [活動名稱] 與您見面！" "期待在[EVENT NAME]见到您！ "
Replace1 "如有任何疑問，請聯繫您的區域經理：" "如有任何疑问，请联系您的区域经理："
Replace1 " [姓名] |  [電子郵件地址] | [WHATSAPP 號碼] (WhatsApp)。 " " [姓名] |  [电子邮件地址] | [WHATSAPP NO] (WhatsApp)。 "
Replace1 "如有任何疑問，請與我們聯繫：" "如有任何疑问，请通过以下方式联系我们："
Replace1 "即時聊天" "实时聊天"
Replace1 "。 " ". "

# -- Comments ---------------------------------------------------------------
# Comment 1 (index 1 / id 0): two paragraphs.
$d.Comments(1).Range.Text = "如果是一天的活动，请选择第一个 `r如果是多天的活动，请选择第二个"
# Comments 2-4 (index 2..4 / ids 1..3): single short paragraph each.
$d.Comments(2).Range.Text = "任选其一"
$d.Comments(3).Range.Text = "任选其一"
$d.Comments(4).Range.Text = "任选其一"
